$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text-valued cells (non-numeric-looking price strings and all percentage cells) ---
$ws.Range("D2").Value = '27.096.61'
$ws.Range("D3").Value = '1.820.32'
$ws.Range("E3").Value = '  -1.50%  '
$ws.Range("E4").Value = '  -1.24%  '
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("E7").Value = '  -1.82%  '
$ws.Range("E8").Value = '  -2.11%  '
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("E10").Value = '  -3.71%  '
$ws.Range("E11").Value = '  -3.55%  '
$ws.Range("D12").Value = '1.832.13'
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("E15").Value = '  -3.25%  '
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("E17").Value = '  -1.28%  '
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").Value = '27.153.68'
$ws.Range("E21").Value = '  -2.07%  '
$ws.Range("E22").Value = '  -2.57%  '
$ws.Range("E23").Value = '  -2.64%  '
$ws.Range("D24").Value = '2.051.40'
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("E26").Value = '  -2.61%  '
$ws.Range("E27").Value = '  +4.95%  '
$ws.Range("E28").Value = '  -1.63%  '
$ws.Range("E30").Value = '  -2.89%  '
$ws.Range("E31").Value = '  -1.92%  '
$ws.Range("E32").Value = '  -4.35%  '
$ws.Range("E33").Value = '  -4.74%  '
$ws.Range("E34").Value = '  +1.56%  '
$ws.Range("E35").Value = '  -3.25%  '
$ws.Range("E36").Value = '  -1.16%  '
$ws.Range("E37").Value = '  -3.56%  '
$ws.Range("E38").Value = '  -0.58%  '
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("E41").Value = '  -0.71%  '
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  -3.16%  '
$ws.Range("E45").Value = '  -1.10%  '
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("E49").Value = '  -1.91%  '
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("E51").Value = '  +1.23%  '

# --- Row 46 / Row 47 swap (Quant <-> Decentraland) ---
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("E47").Value = '  -3.08%  '

# --- Numeric-looking price strings: force text via a scratch cell + PasteSpecial values ---
# (Range.Value auto-converts numeric-looking strings like "1.878" to a Double; routing
#  the write through a Text-formatted helper cell and pasting values keeps the storage
#  type textual, matching the original inlineStr cells, without leaving a cell-level style.)
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"
$helper.Value = '310.46'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Value = '0.4226'
$helper.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$helper.Value = '0.3657'
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Value = '0.07219'
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$helper.Value = '0.8462'
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Value = '20.86'
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Value = '6.648'
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$helper.Value = '0.07079'
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$helper.Value = '5.275'
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$helper.Value = '89.36'
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Value = '1.002'
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$helper.Value = '0.000008826'
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$helper.Value = '14.96'
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$helper.Value = '5.102'
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$helper.Value = '10.82'
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Value = '1.976'
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Value = '151.83'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$helper.Value = '2.247'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$helper.Value = '18.32'
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$helper.Value = '5.212'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$helper.Value = '115.99'
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Value = '0.08791'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$helper.Value = '1.178'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$helper.Value = '0.7415'
$helper.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$helper.Value = '2.966'
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$helper.Value = '4.419'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Value = '1.095'
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$helper.Value = '0.01960'
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$helper.Value = '0.05235'
$helper.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$helper.Value = '7.252'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$helper.Value = '2.865'
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$helper.Value = '0.1688'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Value = '0.5018'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Value = '8.578'
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Value = '10.59'
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$helper.Value = '0.06371'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$helper.Value = '1.655'
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$helper.Value = '1.878'
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$helper.Value = '0.4736'
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$helper.Value = '106.10'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Clear()
$excel.CutCopyMode = $false
